$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.110.63"
$ws.Range("E2").Value = "  +2.00%  "
$ws.Range("D3").Value = "3.103.75"
$ws.Range("E3").Value = "  +4.92%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.98"
$ws.Range("E5").Value = "  +2.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.26"
$ws.Range("E6").Value = "  +5.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.105.26"
$ws.Range("E8").Value = "  +5.12%  "
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("E11").Value = "  +2.70%  "
$ws.Range("E12").Value = "  +6.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.04"
$ws.Range("E14").Value = "  +8.22%  "
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "3.615.99"
$ws.Range("E16").Value = "  +4.88%  "
$ws.Range("D17").Value = "67.174.24"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("E18").Value = "  +3.79%  "
$ws.Range("D19").Value = "3.105.09"
$ws.Range("E19").Value = "  +5.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.06"
$ws.Range("E20").Value = "  +16.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "472.83"
$ws.Range("E21").Value = "  +5.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("E22").Value = "  +6.30%  "
$ws.Range("E23").Value = "  +4.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.79"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  +5.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.85"
$ws.Range("E26").Value = "  +5.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.17"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.11"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.42"
$ws.Range("E30").Value = "  +4.48%  "
$ws.Range("E31").Value = "  +4.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000103"
$ws.Range("E32").Value = "  +5.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.62"
$ws.Range("E33").Value = "  +5.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.116"
$ws.Range("E34").Value = "  +5.74%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  +3.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.92"
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.41"
$ws.Range("E38").Value = "  +9.44%  "
$ws.Range("E39").Value = "  +5.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.32"
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.315"
$ws.Range("E41").Value = "  +4.49%  "
$ws.Range("E42").Value = "  +4.20%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.88"
$ws.Range("E43").Value = "  +2.71%  "
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.72"
$ws.Range("E44").Value = "  +3.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "396.05"
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0364"
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("D47").Value = "2.774.26"
$ws.Range("E47").Value = "  +1.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.37"
$ws.Range("E48").Value = "  +3.67%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.79"
$ws.Range("E50").Value = "  +6.77%  "
$ws.Range("E51").Value = "  +5.21%  "
